$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.259.33"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.368.49"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.693"
$ws.Range("E5").Value = "  +6.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.71"
$ws.Range("E6").Value = "  +3.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.35"
$ws.Range("E7").Value = "  +2.70%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  +28.52%  "

$ws.Range("E10").Value = "  +6.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.89"
$ws.Range("E11").Value = "  +16.96%  "

$ws.Range("E12").Value = "  +19.82%  "

$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.719.97"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "17.12"
$ws.Range("E15").Value = "  +8.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.921"
$ws.Range("E16").Value = "  +7.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.365.55"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.205.44"
$ws.Range("E18").Value = "  +1.66%  "

$ws.Range("E19").Value = "  +4.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.74"
$ws.Range("E20").Value = "  +5.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "78.82"
$ws.Range("E21").Value = "  +5.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "257.78"
$ws.Range("E22").Value = "  +2.64%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("E24").Value = "  +3.49%  "

$ws.Range("E25").Value = "  -2.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.86"
$ws.Range("E26").Value = "  +8.29%  "

$ws.Range("E27").Value = "  +3.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.67"
$ws.Range("E28").Value = "  +9.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.78"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.52"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("E31").Value = "  +2.12%  "

$ws.Range("E32").Value = "  +6.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.42"
$ws.Range("E33").Value = "  +8.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0762"
$ws.Range("E34").Value = "  +10.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.43"
$ws.Range("E35").Value = "  +6.77%  "

$ws.Range("E36").Value = "  +5.78%  "

$ws.Range("E37").Value = "  +1.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.57"
$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("E39").Value = "  +7.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.17"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.12"
$ws.Range("E41").Value = "  +2.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.199"
$ws.Range("E43").Value = "  +18.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.101"
$ws.Range("E44").Value = "  +5.60%  "

$ws.Range("E45").Value = "  +3.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.51"
$ws.Range("E46").Value = "  +12.34%  "

$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.59"
$ws.Range("E48").Value = "  +1.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.48"
$ws.Range("E49").Value = "  -1.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.466.49"
$ws.Range("E50").Value = "  +1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000207"
$ws.Range("E51").Value = "  +2.60%  "
